# Automatische test-sync: 2025-07-23 22:47:50
# Appends a new test-mail log entry (row 25) to the "Logs" sheet and
# updates the corresponding category count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---- 1. Append the new row to the "Logs" worksheet ----
$logs = $wb.Worksheets.Item("Logs")

$row = 25

$logs.Range("A$row").Value = "Ik stuur het pakket morgen terug."
$logs.Range("B$row").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$row").Value = "Testmail #15: Ik stuur het pakket morgen terug."
$logs.Range("D$row").Value = "Retour / Terugbetaling"
$logs.Range("E$row").Value = "Beste klant,`nBedankt voor uw bericht. Om uw retourzending zo soepel mogelijk te laten verlopen, vragen wij u vriendelijk om het volgende te doen:`n- Vul het retourformulier in dat bij uw bestelling zat en voeg dit toe aan het pakket.`n- Stuur het pakket terug naar het volgende adres: [adres retourzending].`n- Zodra wij uw retourzending hebben ontvangen, zullen wij het verder afhandelen en u op de hoogte houden van de status van uw retour.`nMocht u nog verdere vragen of opmerkingen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F$row").Value = "2025-07-23 22:47:16"
$logs.Range("G$row").Value = "Ja"
$logs.Range("H$row").Value = "Nee"
$logs.Range("I$row").Value = "Ja"
$logs.Range("J$row").Value = "Ja"

# The multi-line content in column E triggers an automatic row-height
# bump; AutoFit() brings the row back to the sheet's normal (default)
# height so no stray customHeight is written for the new row.
$logs.Rows.Item($row).AutoFit()

# Extend the conditional-formatting ranges so they keep covering the
# whole data range (D/G/H/I/J columns) now that row 25 was added.
# ModifyAppliesToRange keeps every existing cfRule (and its dxfId) intact
# and simply widens the sqref of the conditionalFormatting block.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$($col)2:$($col)24")
    $newRange = $logs.Range("$($col)2:$($col)25")
    foreach ($fc in $oldRange.FormatConditions) {
        $fc.ModifyAppliesToRange($newRange)
    }
}

# ---- 2. Update the "Dashboard" summary count for "Retour / Terugbetaling" ----
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 4
